$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10:E10").Copy($ws.Range("A11:E11"))

$ws.Range("A11").Value = "33k"
$ws.Range("B11").Value = "R2,R2,R2,R2"
$ws.Range("C11").Value = "R_1206_3216Metric_Pad1.42x1.75mm_HandSolder"
$ws.Range("D11").Value = "C18004"
$ws.Range("E11").Value = 4

$ws.Range("E11").Interior.Pattern = -4142
Write-Host "E11 only done"
